$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "73.232.37"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.10%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.981.52"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.77%  "

# Row 4
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.64"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +5.84%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.29"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +13.00%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.688"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.33%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.793"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.48%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.186"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +8.29%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.46"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +6.34%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000339"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.36%  "

# Row 13
$ws.Range("E13").Value = "  +4.30%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.629.39"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.69%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.983.16"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.80%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.33"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.50%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.26"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.42%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.94"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.22%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.255.67"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.02%  "

# Row 20
$ws.Range("E20").Value = "  -1.10%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "469.42"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.27%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.82"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.39%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "96.53"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.98%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.41"
$ws.Range("D24").ClearFormats()

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.31"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.00%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.26"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.33%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.21"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.64%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.68"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.82%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.95"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.17%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.45"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.17%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.10"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.50%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "14.05"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.19%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "49.81"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.38%  "

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.130"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.09%  "

# Row 35
$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0000103"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +13.76%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "70.58"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.41%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "640.61"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -6.94%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.433"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.33%  "

# Row 39
$ws.Range("E39").Value = "  -0.76%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.41"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.97%  "

# Row 41
$ws.Range("E41").Value = "  +0.20%  "

# Row 42
$ws.Range("E42").Value = "  +0.06%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0487"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.79%  "

# Row 44
$ws.Range("E44").Value = "  +37.98%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.61"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -6.35%  "

# Row 46
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.149"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.75%  "

# Row 47
$ws.Range("B47").Value = "FLOKI"
$ws.Range("C47").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000302"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +9.84%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.43"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.42%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.62"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -5.20%  "

# Row 50
$ws.Range("B50").Value = "WEMIXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.82"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -15.56%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.820.47"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.58%  "
